# Add two new columns ("I0" in I and "IF" in J) to the sheet, matching the
# header formatting used by the existing "IP" column (H), and fill in the
# data rows with the value 8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style/format from H1 into the new header cells so that
# I1/J1 end up sharing the same cell style as the other headers.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))

# Set the new header text.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the data rows.
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8
$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 8
